# Edit slide 4 ("Resource control"), shape 3 ("Content Placeholder 2")
# 1) Remove the two blank paragraphs that used to separate the three
#    existing bullet paragraphs.
# 2) Append two new paragraphs after the "- bumpSwitch_status ..." bullet,
#    the second of which contains several differently-formatted
#    (monospace / colored) runs for inline code snippets.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- smart punctuation -------------------------------------------------
$ldq = [char]0x201C   # “
$rdq = [char]0x201D   # ”

# --- 1) drop the two empty paragraphs (higher index first) -------------
$tr.Paragraphs(4,1).Delete()
$tr.Paragraphs(2,1).Delete()

# --- 2) build the two new paragraphs ------------------------------------
$para4Run1 = "We are using the dcMotor_response which uses the " + $ldq + "status" + $rdq + " returned from "
$para4Run2 = "dcMotor.c"
$para4Run3 = " which is the value returned once the robot has completed its path"
$para4Text = $para4Run1 + $para4Run2 + $para4Run3

$para5Run1  = "There "
$para5Run2  = "are various tasks which have been declared and we use these when we want to perform an operation for example " + $ldq
$para5Run3  = "vTaskSuspend"
$para5Run4  = "("
$para5Run5  = "taskHandle_PlaySong"
$para5Run6  = ")" + $rdq + " and " + $ldq
$para5Run7  = "vTaskResume"
$para5Run8  = "("
$para5Run9  = "taskHandle_PlaySong"
$para5Run10 = ")" + $rdq + " "
$para5Text = $para5Run1 + $para5Run2 + $para5Run3 + $para5Run4 + $para5Run5 + $para5Run6 + $para5Run7 + $para5Run8 + $para5Run9 + $para5Run10

$lastBullet = $tr.Paragraphs(3,1)
$lastBullet.InsertAfter("`r" + $para4Text + "`r" + $para5Text) | Out-Null

# Re-fetch paragraph 4's start position directly -- the InsertAfter return
# value is not reliable for computing offsets in this host.
$base = $tr.Paragraphs(4,1).Start

# --- split paragraph 4 into its three runs ------------------------------
$off = $base
$r = $tr.Characters($off, $para4Run1.Length); $off += $para4Run1.Length
$r = $tr.Characters($off, $para4Run2.Length); $off += $para4Run2.Length
$r = $tr.Characters($off, $para4Run3.Length); $off += $para4Run3.Length

# move past the paragraph mark between paragraph 4 and paragraph 5
$off += 1

# --- split paragraph 5 into its ten runs and format the code runs ------
$r1  = $tr.Characters($off, $para5Run1.Length);  $off += $para5Run1.Length
$r2  = $tr.Characters($off, $para5Run2.Length);  $off += $para5Run2.Length
$r3  = $tr.Characters($off, $para5Run3.Length);  $off += $para5Run3.Length
$r4  = $tr.Characters($off, $para5Run4.Length);  $off += $para5Run4.Length
$r5  = $tr.Characters($off, $para5Run5.Length);  $off += $para5Run5.Length
$r6  = $tr.Characters($off, $para5Run6.Length);  $off += $para5Run6.Length
$r7  = $tr.Characters($off, $para5Run7.Length);  $off += $para5Run7.Length
$r8  = $tr.Characters($off, $para5Run8.Length);  $off += $para5Run8.Length
$r9  = $tr.Characters($off, $para5Run9.Length);  $off += $para5Run9.Length
$r10 = $tr.Characters($off, $para5Run10.Length); $off += $para5Run10.Length

# "vTaskSuspend" ... ")" and "  -> monospace, dark grey, no shadow -> empty effectLst
foreach ($run in @($r3, $r4, $r5, $r6, $r8, $r9, $r10)) {
    $run.Font.Bold = 0
    $run.Font.Italic = 0
    $run.Font.Color.RGB = 0x2F2924
    $run.Font.Shadow = 0
    $run.Font.Name = "ui-monospace"
}

# "vTaskResume" -> monospace, no explicit color, no shadow
$r7.Font.Bold = 0
$r7.Font.Italic = 0
$r7.Font.Shadow = 0
$r7.Font.Name = "ui-monospace"
